# Updates quarterly pressure min/max/avg figures (and the report date)
# for the piezometer table on the active sheet, per "data to update lectures".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M5").Value = 45061
$ws.Range("N14").Value = 52.40547
$ws.Range("O14").Value = 51.76936144390244
$ws.Range("N15").Value = 89.013109
$ws.Range("O15").Value = 88.52062229718875
$ws.Range("O16").Value = 46.22266440526316
$ws.Range("O17").Value = 3.408159194029851
$ws.Range("O18").Value = 54.93904163184079
$ws.Range("O19").Value = -22.21262961764706
$ws.Range("O20").Value = -20.59928452941177
$ws.Range("N21").Value = -6.532022
$ws.Range("O21").Value = -6.872879176470589
$ws.Range("N22").Value = 38.823446
$ws.Range("O22").Value = 38.22670744117647
$ws.Range("N23").Value = 32.441665
$ws.Range("O23").Value = 31.94871207009346
$ws.Range("N24").Value = 39.649686
$ws.Range("O24").Value = 38.40589620175439
$ws.Range("N25").Value = 88.655906
$ws.Range("O25").Value = 88.28820039613527
$ws.Range("O26").Value = -23.6512589112426
$ws.Range("O27").Value = -16.73943071005917
$ws.Range("N28").Value = -7.890887
$ws.Range("O28").Value = -8.288594964497042
$ws.Range("N29").Value = 3.11109
$ws.Range("O29").Value = 2.782040597633136
$ws.Range("N30").Value = 8.42675
$ws.Range("O30").Value = 8.142477834319527
$ws.Range("O31").Value = -5.54305406122449
$ws.Range("O32").Value = -3.929994632653061
$ws.Range("O33").Value = -22.35783700980392
$ws.Range("O34").Value = 56.66089449261084
$ws.Range("N35").Value = -416.56693
$ws.Range("O35").Value = -416.8028627172996
$ws.Range("N36").Value = 152.981757
$ws.Range("O36").Value = 152.675525721519
$ws.Range("O37").Value = -24.69087456540084
$ws.Range("O38").Value = -23.0029835907173
$ws.Range("O39").Value = 60.75640651282051
$ws.Range("N40").Value = -67.69291200000001
$ws.Range("O40").Value = -68.06729568376069
$ws.Range("O41").Value = 22.33319252136752
$ws.Range("N42").Value = 25.445978
$ws.Range("O42").Value = 25.08019196059113
$ws.Range("O43").Value = -9.440129025531915
$ws.Range("N44").Value = -12.247763
$ws.Range("O44").Value = -12.60307751401869
$ws.Range("O45").Value = 43.72242389252337
$ws.Range("O46").Value = -13.8489278
$ws.Range("O47").Value = 25.53825830731707
$ws.Range("O48").Value = -21.77510138514548
$ws.Range("M49").Value = -17.541154
$ws.Range("O49").Value = -17.14029326059214
$ws.Range("O50").Value = -7.437362985294118
$ws.Range("N51").Value = 6.537344
$ws.Range("O51").Value = 6.140783816239316
$ws.Range("O52").Value = -6.689393497872341
$ws.Range("N53").Value = -6.409171
$ws.Range("O53").Value = -6.756637179039301
$ws.Range("O54").Value = -17.07071553586498
$ws.Range("M55").Value = 34.485824
$ws.Range("O55").Value = 34.85136710041841
$ws.Range("O56").Value = 84.13974746443515
$ws.Range("O57").Value = -19.12331364853556
$ws.Range("O58").Value = -99.52137260669456
$ws.Range("O59").Value = 28.46919519323671
$ws.Range("N60").Value = -0.071919
$ws.Range("O60").Value = -1.014941083333333
$ws.Range("O61").Value = -7.271901460639919
$ws.Range("O62").Value = -4.589694835497835
